# Update res_bus vm_pu results for the 380 kV case: the slack/reference
# bus voltage setpoint changed from 1.05 to 1.02 p.u. (column B), which in
# turn shifts the computed per-unit voltages for all other buses
# (columns C:F and I:N) across every time step (rows 2:25).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.048832438469514
$ws.Range("D2").Value = 1.051577634859348
$ws.Range("E2").Value = 1.056448354664799
$ws.Range("F2").Value = 1.068074410413953
$ws.Range("I2").Value = 1.048601033404903
$ws.Range("J2").Value = 1.053874097212401
$ws.Range("K2").Value = 1.054328750161321
$ws.Range("L2").Value = 1.059186037602366
$ws.Range("M2").Value = 1.070780567098123
$ws.Range("N2").Value = 1.021716025764823
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.049889857897529
$ws.Range("D3").Value = 1.052410194190819
$ws.Range("E3").Value = 1.057421644908147
$ws.Range("F3").Value = 1.069203208465161
$ws.Range("I3").Value = 1.048955975781732
$ws.Range("J3").Value = 1.054580023675025
$ws.Range("K3").Value = 1.054973932454646
$ws.Range("L3").Value = 1.059972562450034
$ws.Range("M3").Value = 1.071724495233988
$ws.Range("N3").Value = 1.021958050819042
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.050574120440573
$ws.Range("D4").Value = 1.052948874847748
$ws.Range("E4").Value = 1.058051814650905
$ws.Range("F4").Value = 1.069934280056664
$ws.Range("I4").Value = 1.04918438026571
$ws.Range("J4").Value = 1.055036254886989
$ws.Range("K4").Value = 1.055390719906468
$ws.Range("L4").Value = 1.060481262562683
$ws.Range("M4").Value = 1.072335352584667
$ws.Range("N4").Value = 1.022114305323495
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.050861794570342
$ws.Range("D5").Value = 1.053175325984747
$ws.Range("E5").Value = 1.058316830135869
$ws.Range("F5").Value = 1.070241781398212
$ws.Range("I5").Value = 1.049280098105534
$ws.Range("J5").Value = 1.05522792259861
$ws.Range("K5").Value = 1.055565772291504
$ws.Range("L5").Value = 1.06069506390823
$ws.Range("M5").Value = 1.072592174334372
$ws.Range("N5").Value = 1.022179910326635
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.050910096924928
$ws.Range("D6").Value = 1.053213347516125
$ws.Range("E6").Value = 1.058361332796858
$ws.Range("F6").Value = 1.07029342153948
$ws.Range("I6").Value = 1.049296151757602
$ws.Range("J6").Value = 1.055260096713246
$ws.Range("K6").Value = 1.055595154672277
$ws.Range("L6").Value = 1.060730958796142
$ws.Range("M6").Value = 1.072635296849246
$ws.Range("N6").Value = 1.022190920742363
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.050577964316165
$ws.Range("D7").Value = 1.052951900739569
$ws.Range("E7").Value = 1.058055355439109
$ws.Range("F7").Value = 1.069938388280558
$ws.Range("I7").Value = 1.049185660444454
$ws.Range("J7").Value = 1.055038816478715
$ws.Range("K7").Value = 1.055393059613822
$ws.Range("L7").Value = 1.060484119607662
$ws.Range("M7").Value = 1.072338784181509
$ws.Range("N7").Value = 1.022115182272537
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.049189789867829
$ws.Range("D8").Value = 1.051859010400839
$ws.Range("E8").Value = 1.056777202609114
$ws.Range("F8").Value = 1.068455755523088
$ws.Range("I8").Value = 1.048721250122414
$ws.Range("J8").Value = 1.054112782260087
$ws.Range("K8").Value = 1.054546935004469
$ws.Range("L8").Value = 1.059451895730952
$ws.Range("M8").Value = 1.071099556916639
$ws.Range("N8").Value = 1.021797892037354
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.046743956056497
$ws.Range("D9").Value = 1.049932902848512
$ws.Range("E9").Value = 1.054527900166732
$ws.Range("F9").Value = 1.06584825326342
$ws.Range("I9").Value = 1.047893200124498
$ws.Range("J9").Value = 1.052476781080288
$ws.Range("K9").Value = 1.053050689685976
$ws.Range("L9").Value = 1.057631198968513
$ws.Range("M9").Value = 1.068916440817413
$ws.Range("N9").Value = 1.021236097455618
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.045113588800083
$ws.Range("D10").Value = 1.048648654939387
$ws.Range("E10").Value = 1.053030371529929
$ws.Range("F10").Value = 1.064113338210614
$ws.Range("I10").Value = 1.047334648953143
$ws.Range("J10").Value = 1.051383287094274
$ws.Range("K10").Value = 1.052049661839775
$ws.Range("L10").Value = 1.056416206087551
$ws.Range("M10").Value = 1.067461409293886
$ws.Range("N10").Value = 1.020859767957736
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.044407661661176
$ws.Range("D11").Value = 1.048092523591809
$ws.Range("E11").Value = 1.052382402685791
$ws.Range("F11").Value = 1.06336291101556
$ws.Range("I11").Value = 1.047091244474434
$ws.Range("J11").Value = 1.050909121901566
$ws.Range("K11").Value = 1.051615369527483
$ws.Range("L11").Value = 1.055889817823942
$ws.Range("M11").Value = 1.066831454112936
$ws.Range("N11").Value = 1.020696387678688
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.044145453266433
$ws.Range("D12").Value = 1.047885945340301
$ws.Range("E12").Value = 1.052141788999797
$ws.Range("F12").Value = 1.063084289312423
$ws.Range("I12").Value = 1.047000600573109
$ws.Range("J12").Value = 1.050732894332233
$ws.Range("K12").Value = 1.051453927741822
$ws.Range("L12").Value = 1.055694250412689
$ws.Range("M12").Value = 1.066597473048623
$ws.Range("N12").Value = 1.02063563688025
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.044201697667615
$ws.Range("D13").Value = 1.04793025735819
$ws.Range("E13").Value = 1.052193398248009
$ws.Range("F13").Value = 1.063144049169625
$ws.Range("I13").Value = 1.047020054525898
$ws.Range("J13").Value = 1.050770700338661
$ws.Range("K13").Value = 1.051488563267258
$ws.Range("L13").Value = 1.05573620224984
$ws.Range("M13").Value = 1.06664766222509
$ws.Range("N13").Value = 1.020648671031547
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.04438598734724
$ws.Range("D14").Value = 1.048075447896419
$ws.Range("E14").Value = 1.052362512039039
$ws.Range("F14").Value = 1.063339877604255
$ws.Range("I14").Value = 1.04708375656735
$ws.Range("J14").Value = 1.050894556946835
$ws.Range("K14").Value = 1.051602027275188
$ws.Range("L14").Value = 1.055873653045819
$ws.Range("M14").Value = 1.066812112914274
$ws.Range("N14").Value = 1.020691367305642
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.044499534865326
$ws.Range("D15").Value = 1.048164903750089
$ws.Range("E15").Value = 1.052466718009222
$ws.Range("F15").Value = 1.06346054991046
$ws.Range("I15").Value = 1.047122974668946
$ws.Range("J15").Value = 1.050970855638177
$ws.Range("K15").Value = 1.051671919461843
$ws.Range("L15").Value = 1.055958335268767
$ws.Range("M15").Value = 1.066913438042236
$ws.Range("N15").Value = 1.020717665399412
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.045160439501958
$ws.Range("D16").Value = 1.048685562648106
$ws.Range("E16").Value = 1.053073385013623
$ws.Range("F16").Value = 1.064163158482747
$ws.Range("I16").Value = 1.04735077028497
$ws.Range("J16").Value = 1.051414741646339
$ws.Range("K16").Value = 1.052078466679899
$ws.Range("L16").Value = 1.056451134663538
$ws.Range("M16").Value = 1.067503219082585
$ws.Range("N16").Value = 1.020870601976785
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.045575015885303
$ws.Range("D17").Value = 1.049012146729813
$ws.Range("E17").Value = 1.053454057309336
$ws.Range("F17").Value = 1.064604100985484
$ws.Range("I17").Value = 1.047493245940635
$ws.Range("J17").Value = 1.051692998815975
$ws.Range("K17").Value = 1.052333257953082
$ws.Range("L17").Value = 1.05676017747167
$ws.Range("M17").Value = 1.067873195287127
$ws.Range("N17").Value = 1.020966420756175
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.04581683446037
$ws.Range("D18").Value = 1.049202633483939
$ws.Range("E18").Value = 1.053676142454875
$ws.Range("F18").Value = 1.064861372726832
$ws.Range("I18").Value = 1.047576200210734
$ws.Range("J18").Value = 1.051855236385168
$ws.Range("K18").Value = 1.052481792280581
$ws.Range("L18").Value = 1.056940409097141
$ws.Range("M18").Value = 1.068089004172234
$ws.Range("N18").Value = 1.021022269012366
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.045899288901606
$ws.Range("D19").Value = 1.049267583822394
$ws.Range("E19").Value = 1.053751875455808
$ws.Range("F19").Value = 1.064949108883125
$ws.Range("I19").Value = 1.047604460139036
$ws.Range("J19").Value = 1.05191054415896
$ws.Range("K19").Value = 1.052532424889465
$ws.Range("L19").Value = 1.057001858707443
$ws.Range("M19").Value = 1.068162590790489
$ws.Range("N19").Value = 1.021041304830917
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.045530535437143
$ws.Range("D20").Value = 1.048977107782853
$ws.Range("E20").Value = 1.053413210074607
$ws.Range("F20").Value = 1.064556784019296
$ws.Range("I20").Value = 1.047477975101437
$ws.Range("J20").Value = 1.051663151192281
$ws.Range("K20").Value = 1.052305929650338
$ws.Range("L20").Value = 1.05672702298054
$ws.Range("M20").Value = 1.0678334995127
$ws.Range("N20").Value = 1.020956144570582
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.044331718513313
$ws.Range("D21").Value = 1.048032693073937
$ws.Range("E21").Value = 1.052312710288052
$ws.Range("F21").Value = 1.063282207689886
$ws.Range("I21").Value = 1.047065004326282
$ws.Range("J21").Value = 1.05085808706544
$ws.Range("K21").Value = 1.051568618438808
$ws.Range("L21").Value = 1.055833178410003
$ws.Range("M21").Value = 1.066763685945655
$ws.Range("N21").Value = 1.020678796085587
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.043577999792026
$ws.Range("D22").Value = 1.047438865685604
$ws.Range("E22").Value = 1.051621192004111
$ws.Range("F22").Value = 1.062481527293543
$ws.Range("I22").Value = 1.046804007200366
$ws.Range("J22").Value = 1.050351323982869
$ws.Range("K22").Value = 1.051104310983677
$ws.Range("L22").Value = 1.055270931885084
$ws.Range("M22").Value = 1.066091123014516
$ws.Range("N22").Value = 1.020504045330993
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.043977558147153
$ws.Range("D23").Value = 1.047753668190441
$ws.Range("E23").Value = 1.05198774009024
$ws.Range("F23").Value = 1.062905917180615
$ws.Range("I23").Value = 1.046942494317494
$ws.Range("J23").Value = 1.050620024378508
$ws.Range("K23").Value = 1.051350518401273
$ws.Range("L23").Value = 1.055569013212584
$ws.Range("M23").Value = 1.066447654636149
$ws.Range("N23").Value = 1.020596719128098
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.045550634226952
$ws.Range("D24").Value = 1.048992940387697
$ws.Range("E24").Value = 1.053431667041223
$ws.Range("F24").Value = 1.064578164280132
$ws.Range("I24").Value = 1.047484875797786
$ws.Range("J24").Value = 1.051676638249767
$ws.Range("K24").Value = 1.052318278384531
$ws.Range("L24").Value = 1.056742004153784
$ws.Range("M24").Value = 1.067851436299529
$ws.Range("N24").Value = 1.020960788063845
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.047376228187231
$ws.Range("D25").Value = 1.050430880500788
$ws.Range("E25").Value = 1.055109046019974
$ws.Range("F25").Value = 1.0665217523589
$ws.Range("I25").Value = 1.048108419984956
$ws.Range("J25").Value = 1.052900225121922
$ws.Range("K25").Value = 1.053438128008858
$ws.Range("L25").Value = 1.058102103803433
$ws.Range("M25").Value = 1.06948076199008
$ws.Range("N25").Value = 1.021381652624563
